# Append three new "PBT" broker rows (VCI, 4Q24/1Q25/2Q25) to the bottom of
# Sheet1's data table, then move the view/selection down to the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (which currently ends at row 76).
$newRows = @(
    @("VCI", "4Q24", "PBT", 253.33511597899994),
    @("VCI", "1Q25", "PBT", 355.10702067799997),
    @("VCI", "2Q25", "PBT", 211.49055106)
)

$startRow = 77
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}

# Scroll the view so row 61 is at the top and select the newly added cell A78,
# matching where the author had left the cursor after entering the data.
$win = $excel.ActiveWindow
$win.ScrollRow = 61
$win.ScrollColumn = 1
$ws.Range("A78").Select()
